$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 69
$ws1.Range("F4").Value = 171
$ws1.Range("F6").Value = 5484
$ws1.Range("F7").Value = 121
$ws1.Range("F8").Value = 5423
$ws1.Range("F10").Value = 10
$ws1.Range("F11").Value = 1395
$ws1.Range("F12").Value = 28

# Sheet "全部类型" (4th sheet) - update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 69
$ws4.Range("F4").Value = 171
$ws4.Range("F7").Value = 5484
$ws4.Range("F8").Value = 121
$ws4.Range("F9").Value = 5423
$ws4.Range("F11").Value = 10
$ws4.Range("F12").Value = 1395
$ws4.Range("F13").Value = 28
